# Fixed bug in loads and compile no-fiber loads
# Updates the "conc" (columns B:H) values for rows 4-11 on Sheet1.
# The "conc_raw" values (columns I:O) are left untouched; only the
# normalized/compiled "conc" figures change, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (CCCSD)
$ws.Range("C4").Value = 0.03250201727728265
$ws.Range("D4").Value = 0.01057098089511844
$ws.Range("F4").Value = 0.009783109047814931
$ws.Range("G4").Value = 0.00551796997783287
$ws.Range("H4").Value = 0.002081935221641104

# Row 5 (EBDA)
$ws.Range("B5").Value = 0.006716380654501716
$ws.Range("C5").Value = 0.006011865201232306
$ws.Range("D5").Value = 0.005823994413693795
$ws.Range("F5").Value = 0.008760496505871806
$ws.Range("G5").Value = 0.003317348851699008
$ws.Range("H5").Value = 0.001127224725231057

# Row 6 (EBMUD)
$ws.Range("C6").Value = 0.04984916043789812
$ws.Range("D6").Value = 0.009346932555129099
$ws.Range("F6").Value = 0.01004477436078708
$ws.Range("G6").Value = 0.008996506053690729
$ws.Range("H6").Value = 0.001182682355331754

# Row 7 (FSSD)
$ws.Range("C7").Value = 0.001294559099437148
$ws.Range("D7").Value = 0.000260749703258414
$ws.Range("F7").Value = 0.001468308997471246
$ws.Range("G7").Value = 0.001562209189083905
$ws.Range("H7").Value = 0.0002439024390243902

# Row 8 (PA)
$ws.Range("C8").Value = 0.0004238476641630567
$ws.Range("D8").Value = 0.00009972886215601335
$ws.Range("F8").Value = 0.002576104324475938
$ws.Range("G8").Value = 0.002865857099658951
$ws.Range("H8").Value = 0.00009972886215601333

# Row 9 (SFPUC)
$ws.Range("B9").Value = 0.02802616165560328
$ws.Range("C9").Value = 0.09250247053135975
$ws.Range("D9").Value = 0.0261401247382887
$ws.Range("F9").Value = 0.02015516028062719
$ws.Range("G9").Value = 0.01360349367621611
$ws.Range("H9").Value = 0.0009161436491194176

# Row 10 (SJ)
$ws.Range("B10").Value = 0.004139342847859153
$ws.Range("C10").Value = 0.01097337072326017
$ws.Range("D10").Value = 0.002329613074494138
$ws.Range("E10").Value = 0.0007932220875853444
$ws.Range("F10").Value = 0.002144964288008869
$ws.Range("G10").Value = 0.002241967474901839

# Row 11 (SUNN)
$ws.Range("B11").Value = 0.000992063492063492
$ws.Range("C11").Value = 0.0003306878306878307
$ws.Range("D11").Value = 0.0006613756613756613
$ws.Range("F11").Value = 0.005304101838755305
$ws.Range("G11").Value = 0.005081460527005081
